$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Metadata sheet: Date and FHIR Version values
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# Elements sheet: Extension row (row 2) Constraint(s) text - drop the
# "unless an empty Parameters resource ... or `$this is Parameters" clause
$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Elements sheet: Extension.id row (row 3) Type(s) changes from id to string
$elements.Range("K3").Value = "string`n"

# Elements sheet: Extension.extension row (row 4) Constraint(s) now matches
# the same (deduplicated) ele-1/ext-1 text used in row 2
$elements.Range("AJ4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Elements sheet: Extension.value[x] row (row 6) Definition text switches
# the FHIR version reference in the URL from R4B to R4
$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
